$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select column F (mirrors the user selecting the whole column before
# deleting it) so the saved sheetView selection matches the target.
$ws.Columns("F").Select()

# Drop the per-cell hyperlinks (Customer Email column) before removing
# the column itself.
$ws.Hyperlinks.Delete()

# Bulk import changed Store value for the second data row.
$ws.Range("B3").Value = 5

# Remove the now-unneeded "Customer Email" column entirely.
$ws.Columns("F").Delete()

# Give row 16 a custom height (carried over from the source edit) and
# make sure the sheet's used range/dimension extends down to it without
# leaving a lasting visible formatting change.
$ws.Rows("16").RowHeight = 14.25
$ws.Range("A16").Font.Bold = $true
$ws.Range("A16").Font.Bold = $false
